$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated 2D training schedule values (rows 2-5) and a new row 6.
$data = @(
    @(1, 3, 2, 5, 6, 2, 4, 32, 5),
    @(2, 2, 2, 3, 7, 1, 5, 21, 5),
    @(3, 1, 3, 6, 4, 5, 1, 65, 5),
    @(4, 4, 4, 8, 6, 4, 2, 54, 5),
    @(5, 2, 0, 5, 3, 3, 3, 43, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $vals[$c]
    }
    $ws.Cells.Item($row, 10).Value = "train_dim2_1"
}

# Update the active selection to I1, matching the saved view state.
$ws.Range("I1").Select()
